$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the registration sample data in row 2
$ws.Range("A2").Value = "Kollapudi"
$ws.Range("B2").Value = "Venu"
$ws.Range("C2").Value = "venukollapudi@gmail.com"
$ws.Range("D2").Value = "Venu@12345"
$ws.Range("E2").Value = "Venu@12345"

# Turn the new D2/E2 values into (mailto) hyperlinks, same as the existing C2 one
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Venu@12345")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Venu@12345")

# Re-apply the Hyperlink cell style so D2/E2 look like C2
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("E2").Style = "Hyperlink"

# Move the active selection to G2
$ws.Range("G2").Select() | Out-Null
